$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds values like "29.066.68" / "1.000" / "12.40" that
# Excel would otherwise silently reinterpret as numbers and mangle (dropping
# trailing zeros, re-parsing the dotted groups, etc). Marking the range as
# Text first keeps every write below a faithful, literal string.
$ws.Range("D2:D51").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = '29.066.68'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").Value = '1.835.34'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '242.93'
$ws.Range("E5").Value = '  +0.54%  '
$ws.Range("D6").Value = '0.6266'
$ws.Range("E6").Value = '  -0.69%  '
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.07593'
$ws.Range("E8").Value = '  +3.53%  '
$ws.Range("D9").Value = '0.2932'
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = '22.62'
$ws.Range("E10").Value = '  -1.12%  '
$ws.Range("D11").Value = '0.07743'
$ws.Range("E11").Value = '  +0.83%  '
$ws.Range("D12").Value = '1.843.65'
$ws.Range("E12").Value = '  +0.83%  '
$ws.Range("D13").Value = '4.968'
$ws.Range("E13").Value = '  -0.35%  '
$ws.Range("D14").Value = '0.6655'
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = '82.93'
$ws.Range("E15").Value = '  +1.00%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.000009945'
$ws.Range("E16").Value = '  +14.98%  '
$ws.Range("D17").Value = '6.068'
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("D18").Value = '29.093.41'
$ws.Range("E18").Value = '  +0.25%  '
$ws.Range("D19").Value = '227.22'
$ws.Range("E19").Value = '  +1.42%  '
$ws.Range("D20").Value = '12.40'
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("D22").Value = '7.216'
$ws.Range("E22").Value = '  +1.06%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = '159.47'
$ws.Range("E24").Value = '  +0.91%  '
$ws.Range("D25").Value = '8.515'
$ws.Range("E25").Value = '  +0.75%  '
$ws.Range("D26").Value = '0.1383'
$ws.Range("E26").Value = '  +0.89%  '
$ws.Range("D27").Value = '17.95'
$ws.Range("E27").Value = '  +0.41%  '
$ws.Range("D28").Value = '1.498'
$ws.Range("E28").Value = '  -0.64%  '
$ws.Range("D29").Value = '4.109'
$ws.Range("E29").Value = '  +0.28%  '
$ws.Range("D30").Value = '4.020'
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("D31").Value = '1.194'
$ws.Range("E31").Value = '  -0.67%  '
$ws.Range("D32").Value = '0.05253'
$ws.Range("E32").Value = '  -0.90%  '
$ws.Range("D33").Value = '1.841'
$ws.Range("E33").Value = '  +0.70%  '
$ws.Range("D34").Value = '0.7346'
$ws.Range("E34").Value = '  -0.81%  '
$ws.Range("D35").Value = '1.138'
$ws.Range("E35").Value = '  -1.16%  '
$ws.Range("D36").Value = '2.691'
$ws.Range("E36").Value = '  +1.43%  '
$ws.Range("D37").Value = '1.239.65'
$ws.Range("E37").Value = '  -4.30%  '
$ws.Range("D38").Value = '2.764'
$ws.Range("E38").Value = '  +0.70%  '
$ws.Range("D39").Value = '0.01788'
$ws.Range("D40").Value = '6.373'
$ws.Range("E40").Value = '  +0.60%  '
$ws.Range("D41").Value = '0.8993'
$ws.Range("E41").Value = '  +0.80%  '
$ws.Range("D42").Value = '1.001'
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("D43").Value = '101.94'
$ws.Range("E43").Value = '  -0.60%  '
$ws.Range("D44").Value = '1.985.15'
$ws.Range("E44").Value = '  +0.33%  '
$ws.Range("D45").Value = '0.00000000123'
$ws.Range("E45").Value = '  -0.49%  '
$ws.Range("D46").Value = '64.43'
$ws.Range("E46").Value = '  +0.32%  '
$ws.Range("D47").Value = '0.5114'
$ws.Range("E47").Value = '  -0.49%  '
$ws.Range("D48").Value = '0.4043'
$ws.Range("E48").Value = '  +1.40%  '
$ws.Range("D49").Value = '8.877'
$ws.Range("E49").Value = '  +2.29%  '
$ws.Range("D50").Value = '0.05759'
$ws.Range("E50").Value = '  -1.18%  '
$ws.Range("D51").Value = '6.687'
$ws.Range("E51").Value = '  +0.09%  '
